{"js": "// Update the minimal address line in the contact-info paragraph:\n//   \"herwindo.artono@gmail.com | +628122013050 | Jl. Remaja No.20, Jati Pulogadung, Jakarta Timur\"\n// becomes:\n//   \"herwindo.artono@gmail.com | Jakarta, Indonesia\"\n//\n// The hyperlinked e-mail address (and its run formatting) must stay untouched,\n// so the edit is done as two narrow, unambiguous text replacements that never\n// touch the `w:hyperlink` run.\n\n// 1) Trim \" +628122013050 | Jl. \" down to nothing, leaving just \" | \" in front\n//    of the street address.\nconst prefix = context.document.body.search(\"+628122013050 | Jl. \", {\n  matchCase: true,\n  matchWholeWord: false\n});\nprefix.load(\"text\");\nawait context.sync();\n\nif (prefix.items.length === 0) {\n  throw new Error(\"Could not locate the phone/address prefix to remove.\");\n}\nprefix.items[0].insertText(\"\", Word.InsertLocation.replace);\nawait context.sync();\n\n// 2) Replace the old street address with the new, shorter city/country text.\nconst address = context.document.body.search(\n  \"Remaja No.20, Jati Pulogadung, Jakarta Timur\",\n  { matchCase: true, matchWholeWord: false }\n);\naddress.load(\"text\");\nawait context.sync();\n\nif (address.items.length === 0) {\n  throw new Error(\"Could not locate the old address text to replace.\");\n}\naddress.items[0].insertText(\"Jakarta, Indonesia\", Word.InsertLocation.replace);\nawait context.sync();\n", "ps1": "# Update the minimal address line in the contact-info paragraph:\n#   \"herwindo.artono@gmail.com | +628122013050 | Jl. Remaja No.20, Jati Pulogadung, Jakarta Timur\"\n# becomes:\n#   \"herwindo.artono@gmail.com | Jakarta, Indonesia\"\n#\n# The hyperlinked e-mail address (and its run formatting) must stay untouched,\n# so the edit is done as two narrow, unambiguous Find/Replace passes that\n# never touch the hyperlink run itself.\n\n$d = $word.ActiveDocument\n\n# 1) Remove \" +628122013050 | Jl. \" (the phone number and street prefix),\n#    leaving just \" | \" right after the e-mail hyperlink.\n$find1 = $d.Content\n$find1.Find.Execute(\"+628122013050 | Jl. \", $false, $false, $false, $false, $false, $true, 1, $false, \"\", 2)\n\n# 2) Replace the old street address with the new, shorter city/country text.\n$find2 = $d.Content\n$find2.Find.Execute(\"Remaja No.20, Jati Pulogadung, Jakarta Timur\", $false, $false, $false, $false, $false, $true, 1, $false, \"Jakarta, Indonesia\", 2)\n"}
